# Auto-generated edit script: updates cryptos list values (price/volume)
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range('D2') '96.346.59'
Set-TextValue $ws.Range('E2') '  +4.70%  '
# Row 3
Set-TextValue $ws.Range('D3') '3.598.25'
Set-TextValue $ws.Range('E3') '  +9.70%  '
# Row 4
Set-TextValue $ws.Range('E4') '  -0.06%  '
# Row 5
Set-TextValue $ws.Range('D5') '240.25'
Set-TextValue $ws.Range('E5') '  +6.03%  '
# Row 6
Set-TextValue $ws.Range('D6') '637.86'
Set-TextValue $ws.Range('E6') '  +5.13%  '
# Row 7
Set-TextValue $ws.Range('E7') '  +10.70%  '
# Row 8
Set-TextValue $ws.Range('D8') '0.402'
Set-TextValue $ws.Range('E8') '  +7.22%  '
# Row 9
Set-TextValue $ws.Range('E9') '  -0.06%  '
# Row 10
Set-TextValue $ws.Range('E10') '  +11.23%  '
# Row 11
Set-TextValue $ws.Range('D11') '3.597.55'
Set-TextValue $ws.Range('E11') '  +9.82%  '
# Row 12
Set-TextValue $ws.Range('D12') '43.17'
Set-TextValue $ws.Range('E12') '  +4.93%  '
# Row 13
Set-TextValue $ws.Range('E13') '  +5.44%  '
# Row 14
Set-TextValue $ws.Range('E14') '  +8.83%  '
# Row 15
Set-TextValue $ws.Range('D15') '4.271.60'
Set-TextValue $ws.Range('E15') '  +9.93%  '
# Row 16
Set-TextValue $ws.Range('D16') '96.299.50'
Set-TextValue $ws.Range('E16') '  +4.86%  '
# Row 17
Set-TextValue $ws.Range('D17') '0.0000254'
Set-TextValue $ws.Range('E17') '  +6.37%  '
# Row 18
Set-TextValue $ws.Range('D18') '3.595.42'
Set-TextValue $ws.Range('E18') '  +9.78%  '
# Row 19
Set-TextValue $ws.Range('D19') '13.24'
Set-TextValue $ws.Range('E19') '  +25.60%  '
# Row 20
Set-TextValue $ws.Range('D20') '8.03'
Set-TextValue $ws.Range('E20') '  +1.36%  '
# Row 21
Set-TextValue $ws.Range('D21') '18.11'
Set-TextValue $ws.Range('E21') '  +7.27%  '
# Row 22
Set-TextValue $ws.Range('D22') '0.503'
Set-TextValue $ws.Range('E22') '  +16.16%  '
# Row 23
Set-TextValue $ws.Range('D23') '515.61'
Set-TextValue $ws.Range('E23') '  +7.43%  '
# Row 24
Set-TextValue $ws.Range('D24') '3.46'
Set-TextValue $ws.Range('E24') '  +2.49%  '
# Row 25
Set-TextValue $ws.Range('E25') '  +13.36%  '
# Row 26
Set-TextValue $ws.Range('D26') '6.66'
Set-TextValue $ws.Range('E26') '  +11.60%  '
# Row 27
Set-TextValue $ws.Range('D27') '96.83'
Set-TextValue $ws.Range('E27') '  +9.51%  '
# Row 28
Set-TextValue $ws.Range('D28') '12.46'
Set-TextValue $ws.Range('E28') '  +8.02%  '
# Row 29
Set-TextValue $ws.Range('E29') '  +20.19%  '
# Row 30
Set-TextValue $ws.Range('B30') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D30') '11.60'
Set-TextValue $ws.Range('E30') '  +7.16%  '
# Row 31
Set-TextValue $ws.Range('B31') 'Hedera'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D31') '0.144'
Set-TextValue $ws.Range('E31') '  +6.95%  '
# Row 32
Set-TextValue $ws.Range('E32') '  -0.12%  '
# Row 33
Set-TextValue $ws.Range('E33') '  +7.62%  '
# Row 34
Set-TextValue $ws.Range('D34') '1.00'
Set-TextValue $ws.Range('E34') '  +0.67%  '
# Row 35
Set-TextValue $ws.Range('D35') '30.41'
Set-TextValue $ws.Range('E35') '  +10.19%  '
# Row 36
Set-TextValue $ws.Range('D36') '0.569'
Set-TextValue $ws.Range('E36') '  +10.49%  '
# Row 37
Set-TextValue $ws.Range('D37') '576.53'
Set-TextValue $ws.Range('E37') '  +7.74%  '
# Row 38
Set-TextValue $ws.Range('E38') '  +9.56%  '
# Row 39
Set-TextValue $ws.Range('E39') '  +11.88%  '
# Row 40
Set-TextValue $ws.Range('E40') '  +4.94%  '
# Row 41
Set-TextValue $ws.Range('E41') '  +0.01%  '
# Row 42
Set-TextValue $ws.Range('D42') '0.925'
Set-TextValue $ws.Range('E42') '  +9.93%  '
# Row 43
Set-TextValue $ws.Range('E43') '  +6.74%  '
# Row 44
Set-TextValue $ws.Range('D44') '0.0431'
Set-TextValue $ws.Range('E44') '  +7.74%  '
# Row 45
Set-TextValue $ws.Range('D45') '23.80'
Set-TextValue $ws.Range('E45') '  -0.05%  '
# Row 46
Set-TextValue $ws.Range('D46') '5.66'
Set-TextValue $ws.Range('E46') '  +8.65%  '
# Row 47
Set-TextValue $ws.Range('D47') '3.56'
Set-TextValue $ws.Range('E47') '  +0.29%  '
# Row 48
Set-TextValue $ws.Range('D48') '2.19'
Set-TextValue $ws.Range('E48') '  +7.20%  '
# Row 49
Set-TextValue $ws.Range('D49') '53.84'
Set-TextValue $ws.Range('E49') '  +5.08%  '
# Row 50
Set-TextValue $ws.Range('D50') '8.18'
Set-TextValue $ws.Range('E50') '  +5.38%  '
# Row 51
Set-TextValue $ws.Range('D51') '3.12'
Set-TextValue $ws.Range('E51') '  +6.65%  '
